$d = $word.ActiveDocument

$d.Content.Find.Execute("78-9=69", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=15", 2) | Out-Null
$d.Content.Find.Execute("96-88=8", $true, $false, $false, $false, $false, $true, 1, $false, "3+19=22", 2) | Out-Null
$d.Content.Find.Execute("48-6=42", $true, $false, $false, $false, $false, $true, 1, $false, "72-65=7", 2) | Out-Null
$d.Content.Find.Execute("25+53=78", $true, $false, $false, $false, $false, $true, 1, $false, "19+68=87", 2) | Out-Null
$d.Content.Find.Execute("88-3=85", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=41", 2) | Out-Null
$d.Content.Find.Execute("41-25=16", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("80-33=47", $true, $false, $false, $false, $false, $true, 1, $false, "1+10=11", 2) | Out-Null
$d.Content.Find.Execute("99-40=59", $true, $false, $false, $false, $false, $true, 1, $false, "55+9=64", 2) | Out-Null
$d.Content.Find.Execute("9+49=58", $true, $false, $false, $false, $false, $true, 1, $false, "94-86=8", 2) | Out-Null
$d.Content.Find.Execute("52-47=5", $true, $false, $false, $false, $false, $true, 1, $false, "13+18=31", 2) | Out-Null
$d.Content.Find.Execute("78-34=44", $true, $false, $false, $false, $false, $true, 1, $false, "66-0=66", 2) | Out-Null
$d.Content.Find.Execute("57-14=43", $true, $false, $false, $false, $false, $true, 1, $false, "62-20=42", 2) | Out-Null
$d.Content.Find.Execute("42-0=42", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=88", 2) | Out-Null
$d.Content.Find.Execute("7+9=16", $true, $false, $false, $false, $false, $true, 1, $false, "93-64=29", 2) | Out-Null
$d.Content.Find.Execute("81+8=89", $true, $false, $false, $false, $false, $true, 1, $false, "32+8=40", 2) | Out-Null
$d.Content.Find.Execute("27+19=46", $true, $false, $false, $false, $false, $true, 1, $false, "90-20=70", 2) | Out-Null
$d.Content.Find.Execute("96-25=71", $true, $false, $false, $false, $false, $true, 1, $false, "83-74=9", 2) | Out-Null
$d.Content.Find.Execute("44-8=36", $true, $false, $false, $false, $false, $true, 1, $false, "27+49=76", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $false, $false, $false, $false, $true, 1, $false, "62+33=95", 2) | Out-Null
$d.Content.Find.Execute("94-93=1", $true, $false, $false, $false, $false, $true, 1, $false, "74+25=99", 2) | Out-Null
$d.Content.Find.Execute("18+30=48", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=73", 2) | Out-Null
$d.Content.Find.Execute("65-65=0", $true, $false, $false, $false, $false, $true, 1, $false, "24+50=74", 2) | Out-Null
$d.Content.Find.Execute("78-3=75", $true, $false, $false, $false, $false, $true, 1, $false, "42-11=31", 2) | Out-Null
$d.Content.Find.Execute("42+57=99", $true, $false, $false, $false, $false, $true, 1, $false, "3+62=65", 2) | Out-Null
$d.Content.Find.Execute("68-47=21", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=48", 2) | Out-Null
$d.Content.Find.Execute("85-63=22", $true, $false, $false, $false, $false, $true, 1, $false, "86+8=94", 2) | Out-Null
$d.Content.Find.Execute("15+67=82", $true, $false, $false, $false, $false, $true, 1, $false, "26+53=79", 2) | Out-Null
$d.Content.Find.Execute("49-30=19", $true, $false, $false, $false, $false, $true, 1, $false, "5+32=37", 2) | Out-Null
$d.Content.Find.Execute("11-8=3", $true, $false, $false, $false, $false, $true, 1, $false, "14+48=62", 2) | Out-Null
$d.Content.Find.Execute("98-93=5", $true, $false, $false, $false, $false, $true, 1, $false, "35+49=84", 2) | Out-Null
$d.Content.Find.Execute("60+27=87", $true, $false, $false, $false, $false, $true, 1, $false, "39-9=30", 2) | Out-Null
$d.Content.Find.Execute("39+1=40", $true, $false, $false, $false, $false, $true, 1, $false, "14+0=14", 2) | Out-Null
$d.Content.Find.Execute("33-4=29", $true, $false, $false, $false, $false, $true, 1, $false, "8+55=63", 2) | Out-Null
$d.Content.Find.Execute("69+8=77", $true, $false, $false, $false, $false, $true, 1, $false, "78-42=36", 2) | Out-Null
$d.Content.Find.Execute("33+8=41", $true, $false, $false, $false, $false, $true, 1, $false, "43+26=69", 2) | Out-Null
$d.Content.Find.Execute("74-1=73", $true, $false, $false, $false, $false, $true, 1, $false, "46-44=2", 2) | Out-Null
$d.Content.Find.Execute("48+0=48", $true, $false, $false, $false, $false, $true, 1, $false, "56+15=71", 2) | Out-Null
$d.Content.Find.Execute("13+68=81", $true, $false, $false, $false, $false, $true, 1, $false, "62-46=16", 2) | Out-Null
$d.Content.Find.Execute("17+81=98", $true, $false, $false, $false, $false, $true, 1, $false, "74+20=94", 2) | Out-Null
$d.Content.Find.Execute("65+15=80", $true, $false, $false, $false, $false, $true, 1, $false, "34-32=2", 2) | Out-Null
$d.Content.Find.Execute("59+13=72", $true, $false, $false, $false, $false, $true, 1, $false, "91-69=22", 2) | Out-Null
$d.Content.Find.Execute("31+26=57", $true, $false, $false, $false, $false, $true, 1, $false, "7+87=94", 2) | Out-Null
$d.Content.Find.Execute("56-0=56", $true, $false, $false, $false, $false, $true, 1, $false, "51-18=33", 2) | Out-Null
$d.Content.Find.Execute("64+2=66", $true, $false, $false, $false, $false, $true, 1, $false, "56-4=52", 2) | Out-Null
$d.Content.Find.Execute("28+37=65", $true, $false, $false, $false, $false, $true, 1, $false, "17+65=82", 2) | Out-Null
$d.Content.Find.Execute("65+10=75", $true, $false, $false, $false, $false, $true, 1, $false, "69-47=22", 2) | Out-Null
$d.Content.Find.Execute("91-87=4", $true, $false, $false, $false, $false, $true, 1, $false, "13+62=75", 2) | Out-Null
$d.Content.Find.Execute("42+24=66", $true, $false, $false, $false, $false, $true, 1, $false, "75-12=63", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $false, $false, $false, $false, $true, 1, $false, "76+17=93", 2) | Out-Null
$d.Content.Find.Execute("59-46=13", $true, $false, $false, $false, $false, $true, 1, $false, "13-9=4", 2) | Out-Null
$d.Content.Find.Execute("42+52=94", $true, $false, $false, $false, $false, $true, 1, $false, "67-14=53", 2) | Out-Null
$d.Content.Find.Execute("93-76=17", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=90", 2) | Out-Null
$d.Content.Find.Execute("75-24=51", $true, $false, $false, $false, $false, $true, 1, $false, "82+1=83", 2) | Out-Null
$d.Content.Find.Execute("17+80=97", $true, $false, $false, $false, $false, $true, 1, $false, "43+39=82", 2) | Out-Null
$d.Content.Find.Execute("23-0=23", $true, $false, $false, $false, $false, $true, 1, $false, "95-46=49", 2) | Out-Null
$d.Content.Find.Execute("40+18=58", $true, $false, $false, $false, $false, $true, 1, $false, "32+38=70", 2) | Out-Null
$d.Content.Find.Execute("28+46=74", $true, $false, $false, $false, $false, $true, 1, $false, "15+32=47", 2) | Out-Null
$d.Content.Find.Execute("96-65=31", $true, $false, $false, $false, $false, $true, 1, $false, "58-57=1", 2) | Out-Null
$d.Content.Find.Execute("55-30=25", $true, $false, $false, $false, $false, $true, 1, $false, "17+46=63", 2) | Out-Null
$d.Content.Find.Execute("77+2=79", $true, $false, $false, $false, $false, $true, 1, $false, "48+29=77", 2) | Out-Null
$d.Content.Find.Execute("23+69=92", $true, $false, $false, $false, $false, $true, 1, $false, "94+5=99", 2) | Out-Null
$d.Content.Find.Execute("93-52=41", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("72-29=43", $true, $false, $false, $false, $false, $true, 1, $false, "92+5=97", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "11+11=22", 2) | Out-Null
$d.Content.Find.Execute("24+22=46", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=30", 2) | Out-Null
$d.Content.Find.Execute("63-26=37", $true, $false, $false, $false, $false, $true, 1, $false, "35-28=7", 2) | Out-Null
$d.Content.Find.Execute("65+31=96", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=72", 2) | Out-Null
$d.Content.Find.Execute("43+29=72", $true, $false, $false, $false, $false, $true, 1, $false, "46+19=65", 2) | Out-Null
$d.Content.Find.Execute("89-37=52", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=90", 2) | Out-Null
$d.Content.Find.Execute("94-3=91", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=59", 2) | Out-Null
$d.Content.Find.Execute("3+34=37", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=35", 2) | Out-Null
$d.Content.Find.Execute("83-37=46", $true, $false, $false, $false, $false, $true, 1, $false, "77-64=13", 2) | Out-Null
$d.Content.Find.Execute("64+16=80", $true, $false, $false, $false, $false, $true, 1, $false, "17+30=47", 2) | Out-Null
$d.Content.Find.Execute("35+16=51", $true, $false, $false, $false, $false, $true, 1, $false, "33+28=61", 2) | Out-Null
$d.Content.Find.Execute("29+18=47", $true, $false, $false, $false, $false, $true, 1, $false, "39+53=92", 2) | Out-Null
$d.Content.Find.Execute("85-55=30", $true, $false, $false, $false, $false, $true, 1, $false, "77+20=97", 2) | Out-Null
$d.Content.Find.Execute("84-63=21", $true, $false, $false, $false, $false, $true, 1, $false, "55+43=98", 2) | Out-Null
$d.Content.Find.Execute("21+22=43", $true, $false, $false, $false, $false, $true, 1, $false, "72-17=55", 2) | Out-Null
$d.Content.Find.Execute("47+8=55", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=53", 2) | Out-Null
$d.Content.Find.Execute("74+16=90", $true, $false, $false, $false, $false, $true, 1, $false, "76-59=17", 2) | Out-Null
$d.Content.Find.Execute("47-39=8", $true, $false, $false, $false, $false, $true, 1, $false, "60+33=93", 2) | Out-Null
$d.Content.Find.Execute("56-52=4", $true, $false, $false, $false, $false, $true, 1, $false, "67-16=51", 2) | Out-Null
$d.Content.Find.Execute("24+25=49", $true, $false, $false, $false, $false, $true, 1, $false, "47+20=67", 2) | Out-Null
$d.Content.Find.Execute("44-22=22", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("28+50=78", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=17", 2) | Out-Null
$d.Content.Find.Execute("97-37=60", $true, $false, $false, $false, $false, $true, 1, $false, "71-27=44", 2) | Out-Null
$d.Content.Find.Execute("53-43=10", $true, $false, $false, $false, $false, $true, 1, $false, "43+35=78", 2) | Out-Null
$d.Content.Find.Execute("75-39=36", $true, $false, $false, $false, $false, $true, 1, $false, "77-9=68", 2) | Out-Null
$d.Content.Find.Execute("29+14=43", $true, $false, $false, $false, $false, $true, 1, $false, "46-4=42", 2) | Out-Null
$d.Content.Find.Execute("49-45=4", $true, $false, $false, $false, $false, $true, 1, $false, "73-70=3", 2) | Out-Null
$d.Content.Find.Execute("19+74=93", $true, $false, $false, $false, $false, $true, 1, $false, "22+63=85", 2) | Out-Null
$d.Content.Find.Execute("36-0=36", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("53-29=24", $true, $false, $false, $false, $false, $true, 1, $false, "48-33=15", 2) | Out-Null
$d.Content.Find.Execute("77-70=7", $true, $false, $false, $false, $false, $true, 1, $false, "23+41=64", 2) | Out-Null
$d.Content.Find.Execute("29+30=59", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("29+4=33", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=21", 2) | Out-Null
$d.Content.Find.Execute("17+16=33", $true, $false, $false, $false, $false, $true, 1, $false, "4+11=15", 2) | Out-Null
$d.Content.Find.Execute("8-5=3", $true, $false, $false, $false, $false, $true, 1, $false, "87-78=9", 2) | Out-Null
$d.Content.Find.Execute("61+15=76", $true, $false, $false, $false, $false, $true, 1, $false, "37-32=5", 2) | Out-Null
$d.Content.Find.Execute("97-39=58", $true, $false, $false, $false, $false, $true, 1, $false, "25+50=75", 2) | Out-Null
